# MCH149 collection record: add the data row that was missing from the
# exported sheet (identifier / title / levelOfDescription / extentAndMedium
# / notes), matching the "MCH102 to MCH251" refresh of the archive sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Row 2 values --------------------------------------------------------
$ws.Range("A2").Value = "MCH149"
$ws.Range("C2").Value = "ANC, UMKHONTO WE SIZWE, RE UNION OF POLITICAL PRISONERS NO ROBBEN ISLAND"
$ws.Range("E2").Value = "Series"
$ws.Range("F2").Value = "1 Box"
$ws.Range("G2").Value = "LOCATION: CABINET 1B | GRAP COUNT NUMER: NONE"

# --- Row 2 formatting ------------------------------------------------------
# Body cells use Calibri 10pt, automatic (theme) text colour, no fill -
# matching the header font family but at body size, same as the rest of
# the data rows in this sheet family.
$bodyCells = @("A2", "C2", "D2", "E2", "F2", "G2", "H2")
foreach ($addr in $bodyCells) {
    $cell = $ws.Range($addr)
    $cell.Font.Name = "Calibri"
    $cell.Font.Size = 10
    $cell.Font.ThemeColor = 1
}

# --- Restore the frozen header pane + selection ----------------------------
# (re-selecting + re-freezing keeps the worksheet's split-pane view intact
# after the sheet data changed)
$ws.Range("A2:J2").Select() | Out-Null
$excel.ActiveWindow.FreezePanes = $true
